$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row atom-map labels (columns D,E swap with F,G)
$ws.Range("D1").Value = "N1"
$ws.Range("E1").Value = "N2"
$ws.Range("F1").Value = "C4"
$ws.Range("G1").Value = "C5"

# Fix data rows 2-5 atom-map labels
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "C4"   # C
    $ws.Cells.Item($r, 4).Value = "N5"   # D
    $ws.Cells.Item($r, 5).Value = "N6"   # E
    $ws.Cells.Item($r, 6).Value = "C7"   # F
    $ws.Cells.Item($r, 7).Value = "C16"  # G
    $ws.Cells.Item($r, 8).Value = "C17"  # H
}
